# fix(publipostage): Refactor synthetic array /3
#
# Replace the 4 colored-square emoji used as "statut" markers in column A
# with colored-book emoji, and rename the "noir" (black) status label in
# column B to "bleu" (blue), matching the new "books" icon set:
#   U+1F7E5 (red square)    -> U+1F4D5 (red/closed book)      rouge
#   U+2B1B  (black square)  -> U+1F4D8 (blue book)            noir -> bleu
#   U+1F7E7 (orange square) -> U+1F4D9 (orange book)          orange
#   U+1F7E9 (green square)  -> U+1F4D7 (green book)           vert

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$redSquare    = [string][char]0x1F7E5
$blackSquare  = [string][char]0x2B1B
$orangeSquare = [string][char]0x1F7E7
$greenSquare  = [string][char]0x1F7E9

$redBook    = [string][char]0x1F4D5
$blueBook   = [string][char]0x1F4D8
$orangeBook = [string][char]0x1F4D9
$greenBook  = [string][char]0x1F4D7

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2

    if ($aVal -eq $redSquare) {
        $aCell.Value = $redBook
    } elseif ($aVal -eq $blackSquare) {
        $aCell.Value = $blueBook
    } elseif ($aVal -eq $orangeSquare) {
        $aCell.Value = $orangeBook
    } elseif ($aVal -eq $greenSquare) {
        $aCell.Value = $greenBook
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -eq "noir") {
        $bCell.Value = "bleu"
    }
}
